$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 16, pushing the existing
# rows 16-17 down to 17-18. The new row inherits formatting from the
# row above it (row 15), matching the data rows already in the sheet.
$ws.Rows.Item(16).Insert()

# New row 16 is a fresh record for week 16, with a newer date, using
# the same market/product metadata as its neighboring rows.
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 45173
$ws.Range("D16").NumberFormat = $ws.Range("D17").NumberFormat
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100107
$ws.Range("H16").Value = "Otros"
$ws.Range("I16").Value = 100107002
$ws.Range("J16").Value = "Chirimoya"
$ws.Range("K16").Value = "Cultivar IV Región"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 22000
$ws.Range("P16").Value = 22000
$ws.Range("Q16").Value = "$/bandeja 10 kilos"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 2200
$ws.Range("T16").Value = 10
